$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '57.878.69'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.70%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.454.90'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.11%  '

$ws.Range("E4").Value = '  -0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '511.80'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.79%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '129.95'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.43%  '

$ws.Range("E7").Value = '  -0.19%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.550'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.80%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.469.61'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0962'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.49%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.156'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '5.19'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.33%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.328'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -4.74%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.885.41'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.34%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '57.800.60'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.67%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '21.89'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.44%  '

$ws.Range("E17").Value = '  -2.48%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.457.69'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.52'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.35%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '318.21'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.90%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.12'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("E23").Value = '  +2.40%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '63.01'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.09%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.400'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.87%  '

$ws.Range("E26").Value = '  -0.92%  '

$ws.Range("E27").Value = '  -0.20%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.23'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.86%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '168.89'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.71%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0₃0731'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.69%  '

$ws.Range("E31").Value = '  -2.58%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.17'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.30%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '6.19'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.53%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.993'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.48%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '17.78'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.00%  '

$ws.Range("E37").Value = '  -4.68%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.90'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.56%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '36.56'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("E40").Value = '  -2.14%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.763'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.88%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '271.12'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.21%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.02'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.01%  '

$ws.Range("E44").Value = '  -3.54%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.584'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.34%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0914'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.55%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '120.52'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -5.59%  '

$ws.Range("E48").Value = '  -0.79%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '17.28'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -4.11%  '

$ws.Range("E50").Value = '  -2.10%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '16.66'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.07%  '
